$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "62.785.10"
$ws.Range("E2").Value = "  -1.48%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.679.70"
$ws.Range("E3").Value = "  -1.95%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.00%  "

# Row 5 - BNB
Set-TextValue "D5" "554.40"
$ws.Range("E5").Value = "  -1.41%  "

# Row 6 - Solana
Set-TextValue "D6" "157.00"
$ws.Range("E6").Value = "  -1.78%  "

# Row 7 - USDC
Set-TextValue "D7" "1.00"
$ws.Range("E7").Value = "  +0.02%  "

# Row 8 - XRP
$ws.Range("E8").Value = "  -1.54%  "

# Row 9 - Dogecoin
$ws.Range("E9").Value = "  -3.56%  "

# Row 10 - TRON
$ws.Range("E10").Value = "  -3.32%  "

# Row 12 - Cardano
Set-TextValue "D12" "0.365"
$ws.Range("E12").Value = "  -3.66%  "

# Row 13 - WrappedliquidstakedEther2.0
$ws.Range("D13").Value = "3.153.25"
$ws.Range("E13").Value = "  -1.91%  "

# Row 14 - Avalanche
Set-TextValue "D14" "26.25"
$ws.Range("E14").Value = "  -2.37%  "

# Row 15 - WrappedBTC
$ws.Range("D15").Value = "62.693.08"
$ws.Range("E15").Value = "  -1.35%  "

# Row 16 - ShibaInu
Set-TextValue "D16" "0.0000146"
$ws.Range("E16").Value = "  -2.50%  "

# Row 17 - WrappedEther
$ws.Range("D17").Value = "2.678.86"
$ws.Range("E17").Value = "  -2.10%  "

# Row 18 - Chainlink
Set-TextValue "D18" "11.79"
$ws.Range("E18").Value = "  -5.97%  "

# Row 19 - Polkadot
$ws.Range("E19").Value = "  -3.40%  "

# Row 20 - BitcoinCash
Set-TextValue "D20" "343.86"
$ws.Range("E20").Value = "  -2.99%  "

# Row 21 - Uniswap
$ws.Range("E21").Value = "  -5.66%  "

# Row 22 - Dai
Set-TextValue "D22" "1.00"
$ws.Range("E22").Value = "  +0.02%  "

# Row 23 - Polygon
Set-TextValue "D23" "0.510"
$ws.Range("E23").Value = "  -2.17%  "

# Row 24 - Litecoin
Set-TextValue "D24" "63.05"
$ws.Range("E24").Value = "  -2.07%  "

# Row 25 - Kaspa
$ws.Range("E25").Value = "  -0.50%  "

# Row 26 - Binance-PegBSC-USD
Set-TextValue "D26" "1.00"
$ws.Range("E26").Value = "  +0.07%  "

# Row 27 - InternetComputer(DFINITY)
$ws.Range("E27").Value = "  -3.14%  "

# Row 28 - Fetch.AI
$ws.Range("E28").Value = "  +4.76%  "

# Row 29 - PEPE
$ws.Range("D29").Value = "0.0₃0848"
$ws.Range("E29").Value = "  -6.20%  "

# Row 30 - Aptos
$ws.Range("E30").Value = "  +0.80%  "

# Row 31 - PancakeSwap
$ws.Range("E31").Value = "  -1.49%  "

# Row 32 - Monero
Set-TextValue "D32" "162.46"
$ws.Range("E32").Value = "  -1.92%  "

# Row 33 - USDe
$ws.Range("E33").Value = "  +0.03%  "

# Row 34 - NEARProtocol
Set-TextValue "D34" "4.84"
$ws.Range("E34").Value = "  -1.16%  "

# Row 35 - ImmutableX
$ws.Range("E35").Value = "  -0.86%  "

# Row 36 - EthereumClassic
$ws.Range("E36").Value = "  -3.23%  "

# Row 37 - Stacks
Set-TextValue "D37" "1.78"
$ws.Range("E37").Value = "  -1.01%  "

# Row 38 - Bittensor
Set-TextValue "D38" "339.04"
$ws.Range("E38").Value = "  -2.13%  "

# Row 39 - RenderToken
Set-TextValue "D39" "6.16"
$ws.Range("E39").Value = "  -1.99%  "

# Row 40 - SuiNetwork
$ws.Range("E40").Value = "  -4.81%  "

# Row 41 - Filecoin
$ws.Range("E41").Value = "  -3.17%  "

# Row 42 - OKB
Set-TextValue "D42" "38.33"
$ws.Range("E42").Value = "  -0.57%  "

# Row 43 - InjectiveProtocol
Set-TextValue "D43" "20.77"
$ws.Range("E43").Value = "  -4.99%  "

# Row 44 - now EnergySwap (was Mantle)
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D44" "20.10"
$ws.Range("E44").Value = "  -4.54%  "

# Row 45 - now Mantle (was EnergySwap)
$ws.Range("B45").Value = "Mantle"
$ws.Range("C45").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue "D45" "0.615"
$ws.Range("E45").Value = "  -1.88%  "

# Row 46 - FirstDigitalUSD
$ws.Range("E46").Value = "  -0.07%  "

# Row 47 - Hedera
$ws.Range("E47").Value = "  -5.44%  "

# Row 48 - WhiteBITCoin
Set-TextValue "D48" "11.01"
$ws.Range("E48").Value = "  -0.52%  "

# Row 50 - Aave
Set-TextValue "D50" "128.79"
$ws.Range("E50").Value = "  -2.33%  "

# Row 51 - VeChain
Set-TextValue "D51" "0.0241"
$ws.Range("E51").Value = "  -3.74%  "
